$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) holds plain text like "43.411.91" or "232.51".
# Values that look like a normal decimal number (single '.') would get
# auto-coerced to a numeric cell by plain Value assignment, which would
# not match the source workbook's inlineStr text cells. Force the whole
# column to text first, write the values, then restore the default style
# so the cells end up with no explicit style index (matching the
# original formatting) while still holding string values.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

function Set-Price($row, $value) {
    $ws.Range("D$row").Value = $value
}

function Set-Volume($row, $value) {
    $ws.Range("E$row").Value = $value
}

# Row 2 - Bitcoin
Set-Price 2 "43.373.34"
Set-Volume 2 "  -1.25%  "

# Row 3 - Ethereum
Set-Price 3 "2.364.25"
Set-Volume 3 "  +4.83%  "

# Row 4 - TetherUSD
Set-Volume 4 "  +0.26%  "

# Row 5 - BNB
Set-Price 5 "232.51"
Set-Volume 5 "  +0.31%  "

# Row 6 - XRP
Set-Price 6 "0.649"
Set-Volume 6 "  +0.22%  "

# Row 7 - Solana
Set-Price 7 "68.21"
Set-Volume 7 "  +6.84%  "

# Row 9 - Cardano
Set-Price 9 "0.457"
Set-Volume 9 "  +0.14%  "

# Row 10 - Dogecoin
Set-Price 10 "0.0950"
Set-Volume 10 "  -3.26%  "

# Row 11 - OKB
Set-Price 11 "56.87"
Set-Volume 11 "  +0.12%  "

# Row 12 - Avalanche
Set-Price 12 "26.46"
Set-Volume 12 "  -1.18%  "

# Row 13 - WrappedliquidstakedEther2.0
Set-Price 13 "2.718.08"
Set-Volume 13 "  +4.91%  "

# Row 14 - TRON
Set-Volume 14 "  -0.91%  "

# Row 15 - Chainlink
Set-Price 15 "15.62"
Set-Volume 15 "  +0.19%  "

# Row 16 - Polkadot
Set-Price 16 "6.23"
Set-Volume 16 "  +1.72%  "

# Row 17 - Polygon
Set-Price 17 "0.841"
Set-Volume 17 "  +0.88%  "

# Row 18 - WrappedEther
Set-Price 18 "2.365.23"
Set-Volume 18 "  +4.39%  "

# Row 19 - WrappedBTC
Set-Price 19 "43.385.65"
Set-Volume 19 "  -0.97%  "

# Row 20 - ShibaInu
Set-Price 20 "0.0₃0979"
Set-Volume 20 "  -1.36%  "

# Row 21 - Litecoin
Set-Price 21 "73.99"
Set-Volume 21 "  +0.97%  "

# Row 22 - Uniswap
Set-Price 22 "6.23"
Set-Volume 22 "  +2.81%  "

# Row 23 - BitcoinCash
Set-Price 23 "248.05"
Set-Volume 23 "  -1.14%  "

# Row 24 - WEMIXToken
Set-Volume 24 "  +15.90%  "

# Row 25 - Dai
Set-Volume 25 "  -0.05%  "

# Row 26 - PancakeSwap
Set-Price 26 "2.46"
Set-Volume 26 "  +0.61%  "

# Row 27 - Cosmos
Set-Price 27 "9.95"
Set-Volume 27 "  -0.46%  "

# Row 28 - Toncoin
Set-Volume 28 "  -1.76%  "

# Row 29 - EthereumClassic
Set-Price 29 "22.33"
Set-Volume 29 "  +6.58%  "

# Row 30 - Monero
Set-Price 30 "174.34"
Set-Volume 30 "  +1.97%  "

# Row 31 - ImmutableX
Set-Price 31 "1.53"

# Row 32 - Kaspa
Set-Price 32 "0.128"
Set-Volume 32 "  -6.98%  "

# Row 33 - Stellar
Set-Volume 33 "  -0.24%  "

# Row 34 - Filecoin
Set-Price 34 "4.97"
Set-Volume 34 "  +2.97%  "

# Row 35 - Hedera
Set-Price 35 "0.0690"
Set-Volume 35 "  -2.44%  "

# Row 36 - InternetComputer(DFINITY)
Set-Price 36 "5.05"
Set-Volume 36 "  +2.44%  "

# Row 37 - LidoDAOToken
Set-Price 37 "2.50"
Set-Volume 37 "  +9.12%  "

# Row 38 - THORChain
Set-Price 38 "6.50"
Set-Volume 38 "  +0.90%  "

# Row 39 - RenderToken
Set-Price 39 "3.63"
Set-Volume 39 "  -1.35%  "

# Row 40 - VeChain
Set-Price 40 "0.0255"
Set-Volume 40 "  -2.01%  "

# Row 41 - was FraxShare, now BinanceUSD
$ws.Range("B41").Value = "BinanceUSD"
$ws.Range("C41").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-Price 41 "1.00"
Set-Volume 41 "  -0.04%  "

# Row 42 - was BinanceUSD, now FraxShare
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-Price 42 "8.94"
Set-Volume 42 "  +8.34%  "

# Row 43 - InjectiveProtocol
Set-Price 43 "18.10"
Set-Volume 43 "  +3.89%  "

# Row 44 - ARBITRUM
Set-Volume 44 "  +7.79%  "

# Row 45 - was Aave, now FTXToken
$ws.Range("B45").Value = "FTXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-Price 45 "4.47"
Set-Volume 45 "  +0.76%  "

# Row 46 - was TrustWalletToken, now Aave
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-Price 46 "98.70"
Set-Volume 46 "  +0.97%  "

# Row 47 - was FTXToken, now TrustWalletToken
$ws.Range("B47").Value = "TrustWalletToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-Price 47 "1.21"
Set-Volume 47 "  +1.57%  "

# Row 48 - Cronos
Set-Price 48 "0.0947"
Set-Volume 48 "  +0.39%  "

# Row 49 - Maker
Set-Price 49 "1.442.42"
Set-Volume 49 "  +0.30%  "

# Row 50 - RocketPoolETH
Set-Price 50 "2.590.42"
Set-Volume 50 "  +5.12%  "

# Row 51 - NEARProtocol
Set-Price 51 "2.27"
Set-Volume 51 "  -3.27%  "

# Restore the default (unstyled) look for the price column now that all
# values are written, so cells don't carry a stray explicit style index.
$priceRange.Style = "Normal"
